# Updated symbol list on Sat Jan  7 00:00:11 UTC 2023 with GitHub Actions
#
# This script re-applies the crypto price refresh that was committed: for
# every data row (2-51) the "Data" (F) column moves from 6-1-2023 to
# 7-1-2023 and the "Hora" (G) column resets to 0; the "Price" (D) and
# "Volume(1h)" (E) columns are updated wherever the scraped value actually
# changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $CellRef, $NewValue) {
    $rng = $Worksheet.Range($CellRef)
    # Force a text number format first so Excel does not reinterpret numeric-
    # or date-looking strings (e.g. "259.96", "7-1-2023", "0") as a number or
    # date. Restore the default "Normal" style afterwards so no stray
    # formatting is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

# Per-row overrides for columns D (Price) and E (Volume(1h)).
# Empty string means "leave this column unchanged".
$rowData = @(
    @{ Row = 2;  D = "259.96";     E = "1.30%"   }
    @{ Row = 3;  D = "26.84";      E = "-1.33%"  }
    @{ Row = 4;  D = "4.668";      E = "-0.36%"  }
    @{ Row = 5;  D = "";           E = "3.10%"   }
    @{ Row = 6;  D = "6.695";      E = "0.95%"   }
    @{ Row = 7;  D = "0.8630";     E = "0.57%"   }
    @{ Row = 8;  D = "0.9131";     E = "-2.81%"  }
    @{ Row = 9;  D = "0.1402";     E = "-0.13%"  }
    @{ Row = 10; D = "0.05196";    E = "30.09%"  }
    @{ Row = 11; D = "0.07153";    E = "0.63%"   }
    @{ Row = 12; D = "0.03106";    E = "-1.02%"  }
    @{ Row = 13; D = "0.09147";    E = "0.02%"   }
    @{ Row = 14; D = "0.001534";   E = "0.25%"   }
    @{ Row = 15; D = "0.0006068";  E = "0.45%"   }
    @{ Row = 16; D = "0.006174";   E = "3.19%"   }
    @{ Row = 17; D = "3.486";      E = "-0.91%"  }
    @{ Row = 18; D = "3.180";      E = "-0.69%"  }
    @{ Row = 19; D = "2.185";      E = "-1.29%"  }
    @{ Row = 20; D = "0.3127";     E = "2.46%"   }
    @{ Row = 21; D = "";           E = "-0.16%"  }
    @{ Row = 22; D = "4.107";      E = "7.67%"   }
    @{ Row = 23; D = "0.04257";    E = "0.35%"   }
    @{ Row = 24; D = "";           E = "-0.23%"  }
    @{ Row = 25; D = "0.004037";   E = "-5.87%"  }
    @{ Row = 26; D = "";           E = "0.04%"   }
    @{ Row = 27; D = "";           E = ""        }
    @{ Row = 28; D = "";           E = ""        }
    @{ Row = 29; D = "";           E = ""        }
    @{ Row = 30; D = "";           E = ""        }
    @{ Row = 31; D = "";           E = ""        }
    @{ Row = 32; D = "";           E = ""        }
    @{ Row = 33; D = "";           E = ""        }
    @{ Row = 34; D = "";           E = ""        }
    @{ Row = 35; D = "";           E = ""        }
    @{ Row = 36; D = "";           E = ""        }
    @{ Row = 37; D = "";           E = ""        }
    @{ Row = 38; D = "";           E = ""        }
    @{ Row = 39; D = "";           E = ""        }
    @{ Row = 40; D = "0.03876";    E = "1.34%"   }
    @{ Row = 41; D = "0.1119";     E = "1.53%"   }
    @{ Row = 42; D = "0.004167";   E = "6.26%"   }
    @{ Row = 43; D = "0.01482";    E = "29.50%"  }
    @{ Row = 44; D = "0.002201";   E = "-9.43%"  }
    @{ Row = 45; D = "0.00005336"; E = "-2.38%"  }
    @{ Row = 46; D = "";           E = ""        }
    @{ Row = 47; D = "0.05457";    E = "9.12%"   }
    @{ Row = 48; D = "";           E = "-43.12%" }
    @{ Row = 49; D = "0.00002101"; E = "0.01%"   }
    @{ Row = 50; D = "0.0002001";  E = ""        }
    @{ Row = 51; D = "";           E = ""        }
)

foreach ($entry in $rowData) {
    $row = $entry.Row

    $dRef = "D{0}" -f $row
    $eRef = "E{0}" -f $row
    $fRef = "F{0}" -f $row
    $gRef = "G{0}" -f $row

    if ($entry.D -ne "") {
        Set-TextValue $ws $dRef $entry.D
    }
    if ($entry.E -ne "") {
        Set-TextValue $ws $eRef $entry.E
    }

    # Every row's scrape date moves from 6-1-2023 to 7-1-2023, and the hour
    # resets from 23 to 0.
    Set-TextValue $ws $fRef "7-1-2023"
    Set-TextValue $ws $gRef "0"
}
